# Fix formatting issues introduced when scraping floating point numbers
# and normalize a few stray commas in "Razon social" / "Nombre Fantasia"
# free-text fields (both changes per the source commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Razon social / Nombre Fantasia commas mistakenly used as separators ---
$ws.Range("E6").Value  = "DENING BLANCO. CRISTIAN DAVID"
$ws.Range("F6").Value  = "DENING BLANCO. CRISTIAN DAVID"
$ws.Range("E67").Value  = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E160").Value = "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH"
$ws.Range("E120").Value = "GIMENEZ ANIBAL. FALISTOCCO MARISA DANIELA SH"
$ws.Range("E155").Value = "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO RUBEN"

# --- 2) "Importe" column (H) was scraped with es-AR formatting
#        (1.234,56) and must become plain decimal formatting (1234.56),
#        still stored as text exactly like the source cells. ---
$lastRow = 224
for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 8)
    $old = $cell.Value2
    if ($old -ne $null) {
        $new = $old.Replace(".", "").Replace(",", ".")
        if ($new -ne $old) {
            # Force text storage so "10600.00" isn't re-interpreted as a
            # number (which would drop the trailing zeros / formatting).
            $cell.NumberFormat = "@"
            $cell.Value = $new
            $cell.Style = "Normal"
        }
    }
}
